# "removing agricultural bio gas fired boiler"
#
# HEATING sheet: delete the two district-heating rows for the bio-gas and
# agricultural bio-gas fired boilers (rows 5 and 6). This shifts the
# natural-gas-fired district heating row (old row 7) up into row 5, and lets
# Excel compact the now-unused shared strings ("district heating - bio
# gas-fired boiler", "T23", "district heating - agricultural bio gas-fired
# boiler", "T24") out of sharedStrings.xml automatically.

$wb = $excel.ActiveWorkbook

$wsHeating = $wb.Worksheets.Item("HEATING")
$wsCooling = $wb.Worksheets.Item("COOLING")

$wsHeating.Rows("5:6").Delete()

# Restore the per-sheet cursor positions recorded in the saved view state.
$wsCooling.Activate()
$wsCooling.Range("E20").Select()

$wsHeating.Activate()
$wsHeating.Range("A13").Select()
